$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are numeric-looking text; force text format so
# Excel does not reinterpret/reformat them as numbers (preserving exact string).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.356.38'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.846.84'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.30'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6269'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9987'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07618'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2900'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.79'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07731'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.025'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6793'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001054'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.98'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.141'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.411.03'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.71'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.34'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9986'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.457'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '158.79'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1385'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.430'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.65'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.403'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05602'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.109'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.061'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.162'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.832'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6955'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.586'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01802'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.226.48'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.719'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.357'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9018'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9985'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.28'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.50'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.199'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000117'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.3994'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.983'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.681'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1139'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05699'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4621'

# Volume(1h) percentage strings (column E) already contain surrounding
# whitespace so Excel keeps them as plain text automatically.
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  -1.25%  '
$ws.Range("E9").Value = '  -1.45%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E14").Value = '  -2.85%  '
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  +0.80%  '
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  +7.12%  '
$ws.Range("E28").Value = '  -0.92%  '
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").Value = '  -1.16%  '
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  -2.18%  '
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("E45").Value = '  -4.79%  '
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  -0.14%  '
